$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Program_choosing" sheet used to list several programs (Mannheim Data
# Science, Mannheim Business Informatics, TUM Business Informatics,
# Tuebingen Machine Learning, TU_Berlin_Info_Sys_Mgmt, Uni Muenster
# Information System). This prepares the sheet for Biology/Chemistry: only
# one program entry remains, renamed to "TUM_BIOCHEMISTRY", and the other
# program rows are wiped out (kept as blank rows).

# Rename the remaining program in row 2.
$ws.Cells.Item(2, 1).Value = "TUM_BIOCHEMISTRY"

# Empty out the now unused program rows 3-7 (row stays, contents go away).
$ws.Range("A3:B7").Clear()

# Drop the trailing, never-used rows at the bottom of the sheet (995-1000).
$ws.Range("A995:A1000").EntireRow.Delete()

# The Yes/No dropdown validation on column B only needs to cover the
# remaining two data rows now (was B1:B7).
$ws.Range("B1:B7").Validation.Delete()
$ws.Range("B1:B2").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("B1:B2").Validation.ShowInput = $false
